$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace header "Test" -> "LTSD Parameters" in K2
$ws.Range("K2").Value = "LTSD Parameters"

# New rows of LTSD parameter data (rows shift down: old row5 stays row5,
# but new rows 3 and 4 are inserted with extra parameter headers)
$ws.Range("K3").Value = "Right"
$ws.Range("M3").Value = "Left"

$ws.Range("K4").Value = "Threshols"
$ws.Range("L4").Value = "Win"
$ws.Range("M4").Value = "Threshold"
$ws.Range("N4").Value = "Win"

# Numeric-looking labels must be stored as text (matches the sharedStrings
# entries "8.3", "200.0", "9.0", "280.0"), so force text entry via a
# leading apostrophe and then clear the resulting number format so no
# cell-level style override is left behind.
$ws.Range("L5").Value = "'200.0"
$ws.Range("L5").ClearFormats()
$ws.Range("K5").Value = "'8.3"
$ws.Range("K5").ClearFormats()
$ws.Range("M5").Value = "'9.0"
$ws.Range("M5").ClearFormats()
$ws.Range("N5").Value = "'280.0"
$ws.Range("N5").ClearFormats()

# Update selection to reflect new active cell N5
$ws.Range("N5").Select() | Out-Null
